$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C14").Value = "Ñuble"
$ws.Range("D14").Value = 45079
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E14").Value = 16
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100104
$ws.Range("H14").Value = "Frutos de pepita"
$ws.Range("I14").Value = 100104003
$ws.Range("J14").Value = "Membrillo"
$ws.Range("K14").Value = "Champion"
$ws.Range("L14").Value = "Especial"
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("Q14").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R14").Value = "Región de O'Higgins"
$ws.Range("S14").Value = 667
$ws.Range("T14").Value = 18

# Row 15
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C15").Value = "Ñuble"
$ws.Range("D15").Value = 45079
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100104
$ws.Range("H15").Value = "Frutos de pepita"
$ws.Range("I15").Value = 100104003
$ws.Range("J15").Value = "Membrillo"
$ws.Range("K15").Value = "Champion"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 10000
$ws.Range("Q15").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 556
$ws.Range("T15").Value = 18

# Row 16
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 45079
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100104
$ws.Range("H16").Value = "Frutos de pepita"
$ws.Range("I16").Value = 100104003
$ws.Range("J16").Value = "Membrillo"
$ws.Range("K16").Value = "Champion"
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 20
$ws.Range("N16").Value = 9000
$ws.Range("O16").Value = 9000
$ws.Range("P16").Value = 9000
$ws.Range("Q16").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R16").Value = "Región de O'Higgins"
$ws.Range("S16").Value = 500
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 45069
$ws.Range("D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100104
$ws.Range("H17").Value = "Frutos de pepita"
$ws.Range("I17").Value = 100104003
$ws.Range("J17").Value = "Membrillo"
$ws.Range("K17").Value = "Champion"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("Q17").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R17").Value = "Región de O'Higgins"
$ws.Range("S17").Value = 667
$ws.Range("T17").Value = 18

# Row 18
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C18").Value = "Ñuble"
$ws.Range("D18").Value = 45069
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100104
$ws.Range("H18").Value = "Frutos de pepita"
$ws.Range("I18").Value = 100104003
$ws.Range("J18").Value = "Membrillo"
$ws.Range("K18").Value = "Champion"
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 40
$ws.Range("N18").Value = 10000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 10000
$ws.Range("Q18").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R18").Value = "Región de O'Higgins"
$ws.Range("S18").Value = 556
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C19").Value = "Ñuble"
$ws.Range("D19").Value = 45020
$ws.Range("D19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100104
$ws.Range("H19").Value = "Frutos de pepita"
$ws.Range("I19").Value = 100104003
$ws.Range("J19").Value = "Membrillo"
$ws.Range("K19").Value = "Champion"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 60
$ws.Range("N19").Value = 12000
$ws.Range("O19").Value = 12000
$ws.Range("P19").Value = 12000
$ws.Range("Q19").Value = "`$/caja 18 kilos granel"
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 667
$ws.Range("T19").Value = 18

# Row 20
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 45040
$ws.Range("D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100104
$ws.Range("H20").Value = "Frutos de pepita"
$ws.Range("I20").Value = 100104003
$ws.Range("J20").Value = "Membrillo"
$ws.Range("K20").Value = "Champion"
$ws.Range("L20").Value = "Especial"
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = 13000
$ws.Range("O20").Value = 13000
$ws.Range("P20").Value = 13000
$ws.Range("Q20").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R20").Value = "Región de O'Higgins"
$ws.Range("S20").Value = 722
$ws.Range("T20").Value = 18

# Row 21
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 45040
$ws.Range("D21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100104
$ws.Range("H21").Value = "Frutos de pepita"
$ws.Range("I21").Value = 100104003
$ws.Range("J21").Value = "Membrillo"
$ws.Range("K21").Value = "Champion"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 40
$ws.Range("N21").Value = 12000
$ws.Range("O21").Value = 12000
$ws.Range("P21").Value = 12000
$ws.Range("Q21").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 667
$ws.Range("T21").Value = 18

# Row 22
$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C22").Value = "Ñuble"
$ws.Range("D22").Value = 45070
$ws.Range("D22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100104
$ws.Range("H22").Value = "Frutos de pepita"
$ws.Range("I22").Value = 100104003
$ws.Range("J22").Value = "Membrillo"
$ws.Range("K22").Value = "Champion"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 60
$ws.Range("N22").Value = 10000
$ws.Range("O22").Value = 10000
$ws.Range("P22").Value = 10000
$ws.Range("Q22").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 556
$ws.Range("T22").Value = 18

# Row 23
$ws.Range("A23").Value = 7
$ws.Range("B23").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C23").Value = "Ñuble"
$ws.Range("D23").Value = 45062
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E23").Value = 16
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100104
$ws.Range("H23").Value = "Frutos de pepita"
$ws.Range("I23").Value = 100104003
$ws.Range("J23").Value = "Membrillo"
$ws.Range("K23").Value = "Champion"
$ws.Range("L23").Value = "Especial"
$ws.Range("M23").Value = 50
$ws.Range("N23").Value = 13000
$ws.Range("O23").Value = 13000
$ws.Range("P23").Value = 13000
$ws.Range("Q23").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R23").Value = "Región de O'Higgins"
$ws.Range("S23").Value = 722
$ws.Range("T23").Value = 18

# Row 24
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C24").Value = "Ñuble"
$ws.Range("D24").Value = 45062
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100104
$ws.Range("H24").Value = "Frutos de pepita"
$ws.Range("I24").Value = 100104003
$ws.Range("J24").Value = "Membrillo"
$ws.Range("K24").Value = "Champion"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 12000
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 12000
$ws.Range("Q24").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 667
$ws.Range("T24").Value = 18

# Row 25
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C25").Value = "Ñuble"
$ws.Range("D25").Value = 45033
$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E25").Value = 16
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100104
$ws.Range("H25").Value = "Frutos de pepita"
$ws.Range("I25").Value = 100104003
$ws.Range("J25").Value = "Membrillo"
$ws.Range("K25").Value = "Champion"
$ws.Range("L25").Value = "Especial"
$ws.Range("M25").Value = 60
$ws.Range("N25").Value = 13000
$ws.Range("O25").Value = 13000
$ws.Range("P25").Value = 13000
$ws.Range("Q25").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 722
$ws.Range("T25").Value = 18

# Row 26
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C26").Value = "Ñuble"
$ws.Range("D26").Value = 45033
$ws.Range("D26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E26").Value = 16
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100104
$ws.Range("H26").Value = "Frutos de pepita"
$ws.Range("I26").Value = 100104003
$ws.Range("J26").Value = "Membrillo"
$ws.Range("K26").Value = "Champion"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 80
$ws.Range("N26").Value = 12000
$ws.Range("O26").Value = 12000
$ws.Range("P26").Value = 12000
$ws.Range("Q26").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R26").Value = "Región de O'Higgins"
$ws.Range("S26").Value = 667
$ws.Range("T26").Value = 18

# Row 27
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C27").Value = "Ñuble"
$ws.Range("D27").Value = 45076
$ws.Range("D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100104
$ws.Range("H27").Value = "Frutos de pepita"
$ws.Range("I27").Value = 100104003
$ws.Range("J27").Value = "Membrillo"
$ws.Range("K27").Value = "Champion"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 30
$ws.Range("N27").Value = 12000
$ws.Range("O27").Value = 12000
$ws.Range("P27").Value = 12000
$ws.Range("Q27").Value = "`$/caja 15 kilos granel"
$ws.Range("R27").Value = "Región de O'Higgins"
$ws.Range("S27").Value = 800
$ws.Range("T27").Value = 15

# Row 28
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 45076
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100104
$ws.Range("H28").Value = "Frutos de pepita"
$ws.Range("I28").Value = 100104003
$ws.Range("J28").Value = "Membrillo"
$ws.Range("K28").Value = "Champion"
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 30
$ws.Range("N28").Value = 10000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 10000
$ws.Range("Q28").Value = "`$/caja 15 kilos granel"
$ws.Range("R28").Value = "Región de O'Higgins"
$ws.Range("S28").Value = 667
$ws.Range("T28").Value = 15

# Row 29
$ws.Range("A29").Value = 7
$ws.Range("B29").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C29").Value = "Ñuble"
$ws.Range("D29").Value = 45021
$ws.Range("D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E29").Value = 16
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100104
$ws.Range("H29").Value = "Frutos de pepita"
$ws.Range("I29").Value = 100104003
$ws.Range("J29").Value = "Membrillo"
$ws.Range("K29").Value = "Champion"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 12000
$ws.Range("O29").Value = 12000
$ws.Range("P29").Value = 12000
$ws.Range("Q29").Value = "`$/caja 18 kilos granel"
$ws.Range("R29").Value = "Región de O'Higgins"
$ws.Range("S29").Value = 667
$ws.Range("T29").Value = 18
